$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" footer date from 16.11.17 -> 23.11.17
#    on the slide master and every slide layout (12 placeholders total).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "16.11.17") {
                $tr.Text = "23.11.17"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 1: split "staging.pretrendr.com" into "test" + ".pretrendr.com"
#    (replace the word "staging" with "test" inside the existing run).
# ---------------------------------------------------------------------------
$needle = "staging.pretrendr.com"
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $pos = $full.IndexOf($needle)
        if ($pos -ge 0) {
            # "staging" is the first 7 characters of the needle; Characters()
            # is 1-based, so shift the 0-based IndexOf() result by one.
            $sub = $tr.Characters($pos + 1, 7)
            $sub.Text = "test"
        }
    }
}
